# Natmi following Dr Hou advice:
# Recomputed the LR-pairs (Ccl5-Ccr1) table with an additional "FAPs"
# sending-cluster and an additional "sCs" sending-cluster row-set,
# and refreshed all of the expression/specificity statistics for the
# existing clusters (ECs, M1, M2, Neutro) against the new totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl5"
$ws.Range("C2").Value = "Ccr1"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6472015
$ws.Range("H2").Value = 1.294403
$ws.Range("I2").Value = 0.02102898872844592
$ws.Range("J2").Value = 0.01428114079438904
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 58.378819
$ws.Range("N2").Value = 175.136457
$ws.Range("O2").Value = 0.6920327730022572
$ws.Range("P2").Value = 0.6920327730022573
$ws.Range("Q2").Value = 37.7828592250285
$ws.Range("R2").Value = 226.697155350171
$ws.Range("S2").Value = 0.01455274938317964
$ws.Range("T2").Value = 0.009883017465576707

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl5"
$ws.Range("C3").Value = "Ccr1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6472015
$ws.Range("H3").Value = 1.294403
$ws.Range("I3").Value = 0.02102898872844592
$ws.Range("J3").Value = 0.01428114079438904
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 25.979641
$ws.Range("N3").Value = 77.938923
$ws.Range("O3").Value = 0.3079672269977427
$ws.Range("P3").Value = 0.3079672269977428
$ws.Range("Q3").Value = 16.8140626246615
$ws.Range("R3").Value = 100.884375747969
$ws.Range("S3").Value = 0.006476239345266278
$ws.Range("T3").Value = 0.004398123328812334

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl5"
$ws.Range("C4").Value = "Ccr1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8756340000000001
$ws.Range("H4").Value = 2.626902
$ws.Range("I4").Value = 0.02845125902248994
$ws.Range("J4").Value = 0.02898259453590742
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 58.378819
$ws.Range("N4").Value = 175.136457
$ws.Range("O4").Value = 0.6920327730022572
$ws.Range("P4").Value = 0.6920327730022573
$ws.Range("Q4").Value = 51.11847879624601
$ws.Range("R4").Value = 460.0663091662141
$ws.Range("S4").Value = 0.0196892036767392
$ws.Range("T4").Value = 0.02005690526548408

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl5"
$ws.Range("C5").Value = "Ccr1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.8756340000000001
$ws.Range("H5").Value = 2.626902
$ws.Range("I5").Value = 0.02845125902248994
$ws.Range("J5").Value = 0.02898259453590742
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 25.979641
$ws.Range("N5").Value = 77.938923
$ws.Range("O5").Value = 0.3079672269977427
$ws.Range("P5").Value = 0.3079672269977428
$ws.Range("Q5").Value = 22.74865696739401
$ws.Range("R5").Value = 204.737912706546
$ws.Range("S5").Value = 0.008762055345750733
$ws.Range("T5").Value = 0.008925689270423339

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ccl5"
$ws.Range("C6").Value = "Ccr1"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.540639
$ws.Range("H6").Value = 10.621917
$ws.Range("I6").Value = 0.1150430856889177
$ws.Range("J6").Value = 0.1171915486778959
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 58.378819
$ws.Range("N6").Value = 175.136457
$ws.Range("O6").Value = 0.6920327730022572
$ws.Range("P6").Value = 0.6920327730022573
$ws.Range("Q6").Value = 206.698323325341
$ws.Range("R6").Value = 1860.284909928069
$ws.Range("S6").Value = 0.07961358560403797
$ws.Range("T6").Value = 0.0811003924039933

$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Ccl5"
$ws.Range("C7").Value = "Ccr1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.540639
$ws.Range("H7").Value = 10.621917
$ws.Range("I7").Value = 0.1150430856889177
$ws.Range("J7").Value = 0.1171915486778959
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.979641
$ws.Range("N7").Value = 77.938923
$ws.Range("O7").Value = 0.3079672269977427
$ws.Range("P7").Value = 0.3079672269977428
$ws.Range("Q7").Value = 91.98453013059901
$ws.Range("R7").Value = 827.860771175391
$ws.Range("S7").Value = 0.03542950008487967
$ws.Range("T7").Value = 0.03609115627390258

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Ccl5"
$ws.Range("C8").Value = "Ccr1"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 23.42392733333334
$ws.Range("H8").Value = 70.271782
$ws.Range("I8").Value = 0.761094502822696
$ws.Range("J8").Value = 0.7753081633885378
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 58.378819
$ws.Range("N8").Value = 175.136457
$ws.Range("O8").Value = 0.6920327730022572
$ws.Range("P8").Value = 0.6920327730022573
$ws.Range("Q8").Value = 1367.461214061819
$ws.Range("R8").Value = 12307.15092655637
$ws.Range("S8").Value = 0.5267023393051645
$ws.Range("T8").Value = 0.536538658241057

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Ccl5"
$ws.Range("C9").Value = "Ccr1"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 23.42392733333334
$ws.Range("H9").Value = 70.271782
$ws.Range("I9").Value = 0.761094502822696
$ws.Range("J9").Value = 0.7753081633885378
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 25.979641
$ws.Range("N9").Value = 77.938923
$ws.Range("O9").Value = 0.3079672269977427
$ws.Range("P9").Value = 0.3079672269977428
$ws.Range("Q9").Value = 608.5452229300874
$ws.Range("R9").Value = 5476.907006370786
$ws.Range("S9").Value = 0.2343921635175313
$ws.Range("T9").Value = 0.2387695051474809

$ws.Range("A10").Value = "Neutro"
$ws.Range("B10").Value = "Ccl5"
$ws.Range("C10").Value = "Ccr1"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.243757666666667
$ws.Range("H10").Value = 3.731273
$ws.Range("I10").Value = 0.04041240008444284
$ws.Range("J10").Value = 0.04116711337605242
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 58.378819
$ws.Range("N10").Value = 175.136457
$ws.Range("O10").Value = 0.6920327730022572
$ws.Range("P10").Value = 0.6920327730022573
$ws.Range("Q10").Value = 72.60910370219565
$ws.Range("R10").Value = 653.481933319761
$ws.Range("S10").Value = 0.02796670529411363
$ws.Range("T10").Value = 0.02848899162612787

$ws.Range("A11").Value = "Neutro"
$ws.Range("B11").Value = "Ccl5"
$ws.Range("C11").Value = "Ccr1"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.243757666666667
$ws.Range("H11").Value = 3.731273
$ws.Range("I11").Value = 0.04041240008444284
$ws.Range("J11").Value = 0.04116711337605242
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 25.979641
$ws.Range("N11").Value = 77.938923
$ws.Range("O11").Value = 0.3079672269977427
$ws.Range("P11").Value = 0.3079672269977428
$ws.Range("Q11").Value = 32.31237767099766
$ws.Range("R11").Value = 290.811399038979
$ws.Range("S11").Value = 0.0124456947903292
$ws.Range("T11").Value = 0.01267812174992455

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Ccl5"
$ws.Range("C12").Value = "Ccr1"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.045475
$ws.Range("H12").Value = 2.09095
$ws.Range("I12").Value = 0.03396976365300761
$ws.Range("J12").Value = 0.02306943922721731
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 58.378819
$ws.Range("N12").Value = 175.136457
$ws.Range("O12").Value = 0.6920327730022572
$ws.Range("P12").Value = 0.6920327730022573
$ws.Range("Q12").Value = 61.03359579402501
$ws.Range("R12").Value = 366.2015747641501
$ws.Range("S12").Value = 0.02350818973902214
$ws.Range("T12").Value = 0.01596480800001825

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Ccl5"
$ws.Range("C13").Value = "Ccr1"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.045475
$ws.Range("H13").Value = 2.09095
$ws.Range("I13").Value = 0.03396976365300761
$ws.Range("J13").Value = 0.02306943922721731
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 25.979641
$ws.Range("N13").Value = 77.938923
$ws.Range("O13").Value = 0.3079672269977427
$ws.Range("P13").Value = 0.3079672269977428
$ws.Range("Q13").Value = 27.16106517447501
$ws.Range("R13").Value = 162.96639104685
$ws.Range("S13").Value = 0.01046157391398546
$ws.Range("T13").Value = 0.007104631227199065
